$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New metrics rows for Function-Unit latches / twoToOne mux, appended below
# the existing data (row 152 was the last data row, row 153 is a blank
# separator identical in style to the other section separators in the sheet).
# ---------------------------------------------------------------------------

# Row 153: blank separator row (copy formatting from the very first blank
# separator row in the sheet: row 2, columns B:E).
$ws.Range("B2:E2").Copy() | Out-Null
$ws.Range("B153:E153").PasteSpecial(-4122) | Out-Null

# Write the "instance" labels in the exact order that makes the shared
# string table grow latch1, latch8, latch16, latch32, latch64, twoToOne -
# while still landing each label in its correct row.
$ws.Range("A155").Value = "latch1"
$ws.Range("A156").Value = "latch8"
$ws.Range("A157").Value = "latch16"
$ws.Range("A158").Value = "latch32"
$ws.Range("A159").Value = "latch64"
$ws.Range("A154").Value = "twoToOne"

# Numeric data, column by column.
$ws.Range("B154").Value = 0.636
$ws.Range("B155").Value = 0.545
$ws.Range("B156").Value = 4.388
$ws.Range("B157").Value = 8.716
$ws.Range("B158").Value = 17.461
$ws.Range("B159").Value = 35.029

$ws.Range("C154").Value = 6099.5
$ws.Range("C155").Value = 7.41925
$ws.Range("C156").Value = 64.071699
$ws.Range("C157").Value = 126.268437
$ws.Range("C158").Value = 250.167906
$ws.Range("C159").Value = 544.451882

$ws.Range("D154").Value = 95
$ws.Range("D155").Value = 59
$ws.Range("D156").Value = 59
$ws.Range("D157").Value = 59
$ws.Range("D158").Value = 59
$ws.Range("D159").Value = 59

$ws.Range("E154").Value = 1.85
$ws.Range("E155").Value = 2.117
$ws.Range("E156").Value = 16.935
$ws.Range("E157").Value = 33.87
$ws.Range("E158").Value = 67.74
$ws.Range("E159").Value = 135.481

# ---------------------------------------------------------------------------
# Formatting - replicate the existing styles used by the other "instance"
# rows in the sheet.
# ---------------------------------------------------------------------------

# Column A ("instance" names): style matches A123 / A126 (blue text).
$ws.Range("A123").Copy() | Out-Null
$ws.Range("A154:A159").PasteSpecial(-4122) | Out-Null

# Columns C & D (integer-ish, centered): style matches C127.
$ws.Range("C127").Copy() | Out-Null
$ws.Range("C154:C159").PasteSpecial(-4122) | Out-Null
$ws.Range("D154:D159").PasteSpecial(-4122) | Out-Null

# Columns B & E (two decimal places, centered): start from the same base
# style as column C/D, then switch the number format to 0.00 to match the
# new style added to the workbook.
$ws.Range("B154:B159").PasteSpecial(-4122) | Out-Null
$ws.Range("E154:E159").PasteSpecial(-4122) | Out-Null
$ws.Range("B154:B159").NumberFormat = "0.00"
$ws.Range("E154:E159").NumberFormat = "0.00"

# Keep the selection / active cell roughly where the author left it after
# appending the new rows.
$ws.Range("A168").Select() | Out-Null
